$p = $ppt.ActivePresentation

# --- Slide 3: "Interface Progress" -> Content Placeholder 5 (shape id 6, collection index 5) ---
$s3 = $p.Slides.Item(3)
$shape3 = $s3.Shapes.Item(5)
$tf3 = $shape3.TextFrame.TextRange

# Rebuild the paragraphs for this text box.
$tf3.Text = "Major problems running current build on different computers`rMade as a website using ASP.net`rPages built out for login, inventory adding, account registration, viewing inventory and a min menu `r`r"

# --- Slide 5: "Scanner Progress" -> Content Placeholder 5 (shape id 6, collection index 5) ---
$s5 = $p.Slides.Item(5)
$shape5 = $s5.Shapes.Item(5)
$tf5 = $shape5.TextFrame.TextRange
$tf5.Text = "Candidate APIs found but do to problems running the interface have not been implamented"
